$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45 (shifts old rows 45-49 down to 46-50),
# copying the row's cell formatting from row 33 (same "middle of merged
# group, with E/F/G/H/I narrower borders" visual pattern) so the new row
# keeps the look of its siblings.
$ws.Rows.Item(45).Insert(-4121, 0)

$ws.Range("A33:J33").Copy()
$ws.Range("A45:J45").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's content.
$ws.Range("C45").Value = 43
$ws.Range("D45").Value = "Cliente não autoriza operar animal"
$ws.Range("E45").Value = "X(37)"

# Renumber the sequence column (C) for the rows that shifted down, since
# Excel's row-insert keeps their original literal numbers.
$ws.Range("C46").Value = 44
$ws.Range("C47").Value = 45
$ws.Range("C48").Value = 46
$ws.Range("C49").Value = 47
$ws.Range("C50").Value = 48

$ws.Range("C48").Select()
